$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ch 1")

# Insert a new blank row at position 7, pushing rows 7-20 down to 8-21
$ws.Rows.Item(7).Insert()

# Update row 2 (previously "data_sci_neurosci" entry, now "why" entry with more data)
$ws.Range("B2").Value = "why"
$ws.Range("C2").Value = 363
$ws.Range("D2").Value = "20 min"
$ws.Range("D2").HorizontalAlignment = -4108  # xlCenter
$ws.Range("E2").Value = "y"
$ws.Range("F2").Value = "?"

# Update rows 3-5
$ws.Range("B3").Value = "learning objectives (of course)"
$ws.Range("B4").Value = "philosophy"
$ws.Range("B5").Value = "Practical skills"

# Update rows 8-15 (shifted down by the inserted row)
$ws.Range("B8").Value = "teaching_approach"
$ws.Range("B9").Value = "online_learning"
$ws.Range("B10").Value = "constructivism, connectionism"
$ws.Range("B11").Value = "core_principles"
$ws.Range("B12").Value = "values_goals"
$ws.Range("B13").Value = "learning"
$ws.Range("B14").Value = "mindset"
$ws.Range("B15").Value = "erros_debugging"

# Row 16 is blank; row 17 gets "teamwork"
$ws.Range("B17").Value = "teamwork"

# Fix row height for the wrap-text row (now row 10)
$ws.Rows.Item(10).RowHeight = 17

# The "new session" date-style marker (direct format s=7) moves from row 13 back to row 12
$ws.Range("A13").Copy() | Out-Null
$ws.Range("A12").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("A6").Copy() | Out-Null
$ws.Range("A13").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Row 20 should end up completely blank (no formatting carried over from the insert)
$ws.Range("A20:F20").Clear() | Out-Null

# Fix the SUM formula to start from C3 instead of C2 (row insert already pointed it at C2:C20)
$ws.Range("C21").Formula = "=SUM(C3:C20)"

# Update selection
$ws.Range("B5").Select()
